$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new "season record" columns (AD, AE, AF),
# matching the bold/centered/bordered style used by the other headers in row 1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 45; $r++) {
    $ws.Range("AD$r").Value = 77
    $ws.Range("AE$r").Value = 85
    $ws.Range("AF$r").Value = 0
}

Write-Output "done"
